# ---------------------------------------------------------------------------
# Applies the "Continue with GATech Modeling, UCSD Python for DS, Udacity
# Inferential Stats, and Udemy SQL Courses" commit to Lesson8_ANOVA_Continued
# workbook:
#   * adds ANOVA post-hoc (Tukey HSD) calculations to the "food" sheet
#   * adds a new "cancer" sheet (one-way ANOVA on drug/placebo data)
#   * adds a new "quiz" sheet (F-crit lookup quiz)
#   * makes "cancer" the active sheet/tab
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$food = $wb.Worksheets.Item(1)
$food.Name = "food"

# A cell that already carries cellXfs style index 1 (bold label, no border) -
# used as a formatting template for the new bold label cells we add below.
$boldTemplate = $food.Range("E2")

function Set-BoldLabel($range, $value) {
    $boldTemplate.Copy() | Out-Null
    $range.PasteSpecial(-4122) | Out-Null
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# food sheet: Tukey HSD post-hoc test + eta^2 effect size
# ---------------------------------------------------------------------------

# S(p) = SUM(I5)  (pooled/average within-group variance, already computed above)
Set-BoldLabel $food.Range("I14") "S(p)"
$food.Range("J14").Copy() | Out-Null
$boldTemplate.Copy() | Out-Null
$food.Range("J14").PasteSpecial(-4122) | Out-Null
$food.Range("J14").Formula = "=SUM(I5)"

# q (from tbl) and the note about table lookup
Set-BoldLabel $food.Range("I15") "q (from tbl)"
$boldTemplate.Copy() | Out-Null
$food.Range("J15").PasteSpecial(-4122) | Out-Null
$food.Range("J15").Value = 4.34
$food.Range("K15").Value = '(# of groups and "dF(w)" or "dF(error term)")'

# HSD = q * SQRT( S(p) / n )
Set-BoldLabel $food.Range("I16") "HSD"
$boldTemplate.Copy() | Out-Null
$food.Range("J16").PasteSpecial(-4122) | Out-Null
$food.Range("J16").Formula = "=J15*SQRT(I5/F3)"

$food.Range("J17").Value = "* If any 2 samples have a mean difference greater than HSD, the difference is honestly significant"

$boldTemplate.Copy() | Out-Null
$food.Range("J18").PasteSpecial(-4122) | Out-Null
$food.Range("J18").Value = "' '- 3, 6, 9 --> all honestly significantly different --> most food eaten w/ C --> cows prefer food C"

$food.Range("J19").Value = "* CAN ONLY COMPUTE WHEN ALL SAMPLE SIZES ARE THE SAME"

# pairwise mean-difference comparisons (Cd = comparison difference, in HSD units)
Set-BoldLabel $food.Range("I21") "Cd(ab)"
$boldTemplate.Copy() | Out-Null
$food.Range("J21").PasteSpecial(-4122) | Out-Null
$food.Range("J21").Formula = "=(A5-B5)/`$J`$14"

Set-BoldLabel $food.Range("I22") "Cd(bc)"
$boldTemplate.Copy() | Out-Null
$food.Range("J22").PasteSpecial(-4122) | Out-Null
$food.Range("J22").Formula = "=(A5-C5)/`$J`$14"

Set-BoldLabel $food.Range("I23") "Cd(ac)"
$boldTemplate.Copy() | Out-Null
$food.Range("J23").PasteSpecial(-4122) | Out-Null
$food.Range("J23").Formula = "=(B5-C5)/`$J`$14"

# eta^2 effect size = SS(b) / (SS(b)+SS(w))
Set-BoldLabel $food.Range("F25") "eta^2"
$food.Range("G25").Formula = "=F4/(F4+F5)"

Set-BoldLabel $food.Range("G26") "* 90% of total variation in amount eaten is due to BG differences/differences in food types, and 10% remains unexplained"
Set-BoldLabel $food.Range("G27") "* anything bigger than 0.4 is pretty big"

$food.Range("G25").Select() | Out-Null

# ---------------------------------------------------------------------------
# new "cancer" sheet: one-way ANOVA of tumor-reduction by drug type
# ---------------------------------------------------------------------------

$cancer = $wb.Worksheets.Add($null, $food)
$cancer.Name = "cancer"

$cancer.Range("A1").Value = "placebo"
$cancer.Range("B1").Value = "drug1"
$cancer.Range("C1").Value = "drug2"
$cancer.Range("D1").Value = "drug3"

$cancer.Range("A2").Value = 1.5
$cancer.Range("A3").Value = 1.3
$cancer.Range("A4").Value = 1.8
$cancer.Range("A5").Value = 1.6
$cancer.Range("A6").Value = 1.3

$cancer.Range("B2").Value = 1.6
$cancer.Range("B3").Value = 1.7
$cancer.Range("B4").Value = 1.9
$cancer.Range("B5").Value = 1.2

$cancer.Range("C2").Value = 2
$cancer.Range("C3").Value = 1.4
$cancer.Range("C4").Value = 1.5
$cancer.Range("C5").Value = 1.5
$cancer.Range("C6").Value = 1.8
$cancer.Range("C7").Value = 1.7
$cancer.Range("C8").Value = 1.4

$cancer.Range("D2").Value = 2.9
$cancer.Range("D3").Value = 3.1
$cancer.Range("D4").Value = 2.8
$cancer.Range("D5").Value = 2.7

$boldTemplate.Copy() | Out-Null
$cancer.Range("F2").PasteSpecial(-4122) | Out-Null
$cancer.Range("F2").Value = "mu(g)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G2").PasteSpecial(-4122) | Out-Null
$cancer.Range("G2").Formula = "=AVERAGE(A2:D8)"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F3").PasteSpecial(-4122) | Out-Null
$cancer.Range("F3").Value = "SS(b)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G3").PasteSpecial(-4122) | Out-Null
$cancer.Range("G3").Formula = "=SUM(A16:D16)"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F4").PasteSpecial(-4122) | Out-Null
$cancer.Range("F4").Value = "SS(w)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G4").PasteSpecial(-4122) | Out-Null
$cancer.Range("G4").Formula = "=SUM(A28:D34)"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F5").PasteSpecial(-4122) | Out-Null
$cancer.Range("F5").Value = "dF(b)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G5").PasteSpecial(-4122) | Out-Null
$cancer.Range("G5").Formula = "=COUNT(A9:D9)-1"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F6").PasteSpecial(-4122) | Out-Null
$cancer.Range("F6").Value = "dF(w)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G6").PasteSpecial(-4122) | Out-Null
$cancer.Range("G6").Formula = "=COUNT(A2:D8)-(G5+1)"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F7").PasteSpecial(-4122) | Out-Null
$cancer.Range("F7").Value = "MS(b)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G7").PasteSpecial(-4122) | Out-Null
$cancer.Range("G7").Formula = "=G3/G5"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F8").PasteSpecial(-4122) | Out-Null
$cancer.Range("F8").Value = "MS(w)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G8").PasteSpecial(-4122) | Out-Null
$cancer.Range("G8").Formula = "=G4/G6"

$boldTemplate.Copy() | Out-Null
$cancer.Range("A9:D9").PasteSpecial(-4122) | Out-Null
$cancer.Range("A9").Formula = "=AVERAGE(A2:A8)"
$cancer.Range("B9:D9").Formula = "=AVERAGE(B2:B8)"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F9").PasteSpecial(-4122) | Out-Null
$cancer.Range("F9").Value = "F"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G9").PasteSpecial(-4122) | Out-Null
$cancer.Range("G9").Formula = "=G7/G8"

$boldTemplate.Copy() | Out-Null
$cancer.Range("F10").PasteSpecial(-4122) | Out-Null
$cancer.Range("F10").Value = "eta^2"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G10").PasteSpecial(-4122) | Out-Null
$cancer.Range("G10").Formula = "=G3/SUM(G3:G4)"

$cancer.Range("A11").Value = "var(b)"
$boldTemplate.Copy() | Out-Null
$cancer.Range("G11").PasteSpecial(-4122) | Out-Null
$cancer.Range("G11").Value = "* 86% of differences in tumor reductions/dependent variable is due to the type of drug given/can be explained by the different treatments"

$boldTemplate.Copy() | Out-Null
$cancer.Range("A12:D12").PasteSpecial(-4122) | Out-Null
$cancer.Range("A12").Formula = "=A9-`$G`$2"
$cancer.Range("B12:D12").Formula = "=B9-`$G`$2"

$cancer.Range("A13").Value = "var(b)^2"

$boldTemplate.Copy() | Out-Null
$cancer.Range("A14:D14").PasteSpecial(-4122) | Out-Null
$cancer.Range("A14").Formula = "=A12^2"
$cancer.Range("B14").Formula = "=B12^2"
$cancer.Range("C14").Formula = "=C12^2"
$cancer.Range("D14").Formula = "=D12^2"

$cancer.Range("A15").Value = "var(b)^2 * n"

$cancer.Range("A16").Formula = "=A14*COUNT(A2:A8)"
$cancer.Range("B16:C16").Formula = "=B14*COUNT(B2:B8)"
$cancer.Range("D16").Formula = "=D14*COUNT(D2:D8)"

$cancer.Range("A18").Value = "var(w)"

$cancer.Range("A19").Formula = "=A2-`$A`$9"
$cancer.Range("B19").Formula = "=B2-`$B`$9"
$cancer.Range("C19").Formula = "=C2-`$C`$9"
$cancer.Range("D19").Formula = "=D2-`$D`$9"

$cancer.Range("A20:A24").Formula = "=A3-`$A`$9"
$cancer.Range("B20:B22").Formula = "=B3-`$B`$9"
$cancer.Range("C20:C25").Formula = "=C3-`$C`$9"
$cancer.Range("D20:D22").Formula = "=D3-`$D`$9"

$cancer.Range("C25").Formula = "=C8-`$C`$9"

$cancer.Range("A27").Value = "var(w)^2"
$boldTemplate.Copy() | Out-Null
$cancer.Range("B27").PasteSpecial(-4122) | Out-Null

$grayTemplate = $food.Range("E11")
$grayTemplate.Copy() | Out-Null
$cancer.Range("A28:D28").PasteSpecial(-4122) | Out-Null
$cancer.Range("A28").Formula = "=A19^2"
$cancer.Range("B28:D28").Formula = "=B19^2"

$grayTemplate.Copy() | Out-Null
$cancer.Range("A29:D29").PasteSpecial(-4122) | Out-Null
$cancer.Range("A29:D29").Formula = "=A20^2"

$grayTemplate.Copy() | Out-Null
$cancer.Range("A30:D30").PasteSpecial(-4122) | Out-Null
$cancer.Range("A30:D30").Formula = "=A21^2"

$grayTemplate.Copy() | Out-Null
$cancer.Range("A31:D31").PasteSpecial(-4122) | Out-Null
$cancer.Range("A31:D31").Formula = "=A22^2"

$grayTemplate.Copy() | Out-Null
$cancer.Range("A32:D32").PasteSpecial(-4122) | Out-Null
$cancer.Range("A32").Formula = "=A23^2"
$cancer.Range("C32").Formula = "=C23^2"

$grayTemplate.Copy() | Out-Null
$cancer.Range("A33:D33").PasteSpecial(-4122) | Out-Null
$cancer.Range("C33").Formula = "=C24^2"

$grayTemplate.Copy() | Out-Null
$cancer.Range("A34:D34").PasteSpecial(-4122) | Out-Null
$cancer.Range("C34").Formula = "=C25^2"

$cancer.Columns("A").AutoFit() | Out-Null

$cancer.Range("G12").Select() | Out-Null

# ---------------------------------------------------------------------------
# new "quiz" sheet: F-crit table lookup quiz (two worked examples)
# ---------------------------------------------------------------------------

$quiz = $wb.Worksheets.Add($null, $cancer)
$quiz.Name = "quiz"

$quiz.Range("A1").Value = "df(b)"
$quiz.Range("B1").Value = 2
$quiz.Range("D1").Value = "df(b)"
$quiz.Range("E1").Value = 3

$quiz.Range("A2").Value = "df(w)"
$quiz.Range("B2").Value = 30
$quiz.Range("D2").Value = "df(w)"
$quiz.Range("E2").Value = 15

$quiz.Range("A3").Value = "k"
$quiz.Range("B3").Formula = "=B1+1"
$quiz.Range("D3").Value = "k"
$quiz.Range("E3").Formula = "=E1+1"

$quiz.Range("A4").Value = "alpha  "
$quiz.Range("B4").Value = 0.05
$quiz.Range("D4").Value = "alpha  "
$quiz.Range("E4").Value = 0.05

$boldTemplate.Copy() | Out-Null
$quiz.Range("A5:B5").PasteSpecial(-4122) | Out-Null
$quiz.Range("A5").Value = "F-crit"
$quiz.Range("B5").Value = 2.89

$boldTemplate.Copy() | Out-Null
$quiz.Range("D5:E5").PasteSpecial(-4122) | Out-Null
$quiz.Range("D5").Value = "F-crit"
$quiz.Range("E5").Value = 2.89

$quiz.Range("E3").Select() | Out-Null

$cancer.Activate() | Out-Null
